# Insert a new record row at row 25 (pushes existing rows 25..109 down to 26..110)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(25).Insert()

# Populate the newly inserted row 25 with the new Granada price record
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 44659
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100104
$ws.Range("H25").Value = "Frutos de pepita"
$ws.Range("I25").Value = 100104001
$ws.Range("J25").Value = "Granada"
$ws.Range("K25").Value = "Wonderfull"
$ws.Range("L25").Value = "Especial"
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 25000
$ws.Range("O25").Value = 25000
$ws.Range("P25").Value = 25000
$ws.Range("Q25").Value = '$/caja 18 kilos empedrada'
$ws.Range("R25").Value = "Provincia de Limarí"
$ws.Range("S25").Value = 1389
$ws.Range("T25").Value = 18
